$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "model" (sheet1): item_pack_barcode field's type changes from
# "string" to "number".
# ---------------------------------------------------------------------------
$wsModel = $wb.Worksheets.Item("model")
$wsModel.Range("A7").Value = "number"

# ---------------------------------------------------------------------------
# Sheet "survey" (sheet3): a new "confirm/update barcode" prompt is inserted
# right before the old "item pack barcode out of range" branch-check row
# (old row 18), pushing everything from the old row 18 onward down by one.
# ---------------------------------------------------------------------------
$wsSurvey = $wb.Worksheets.Item("survey")

# Insert a new blank row at row 18 (old row 18 and below shift to 19+).
$wsSurvey.Range("A18").EntireRow.Insert()

# The freshly inserted row inherits formatting/content from the row above it
# (row 17) including a stray H18 cell - clear that first, then re-apply the
# formatting used by the equivalent "note" rows (e.g. row 16: string/
# item_pack_barcode/label pattern in columns C/E/F).
$wsSurvey.Range("H18").Clear()

$wsSurvey.Range("C16").Copy()
$wsSurvey.Range("C18").PasteSpecial(-4122)
$wsSurvey.Range("E16").Copy()
$wsSurvey.Range("E18").PasteSpecial(-4122)
$wsSurvey.Range("F16").Copy()
$wsSurvey.Range("F18").PasteSpecial(-4122)

# New row 18 content: a "confirm or update" note for item_pack_barcode.
$wsSurvey.Range("C18").Value = "string"
$wsSurvey.Range("E18").Value = "item_pack_barcode"
$wsSurvey.Range("F18").Value = "Please confirm or update item pack barcode"

# Old row 18 (now row 19) keeps its "if" condition shape, but the condition
# text itself is rewritten to use the new field/range names.
$wsSurvey.Range("B19").Value = "(data('item_pack_barcode') > data('max_range')) || (data('item_pack_barcode') < data('min_range'))"

# ---------------------------------------------------------------------------
# View/selection state: "model" was the active tab with A16 selected; now
# "survey" is the active tab (tabSelected moves there) with B19 selected,
# and "model" is left with B15 selected.
# ---------------------------------------------------------------------------
$wsModel.Activate()
$wsModel.Range("B15").Select()

$wsSurvey.Activate()
$wsSurvey.Range("B19").Select()
